$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WHO")

# New rows to append (WHO_code, WHO_countryname, wbcountryname, wbcode)
# ASM / American Samoa
$ws.Range("A207").Value = "ASM"
$ws.Range("B207").Value = "American Samoa"
$ws.Range("C207").Value = "American Samoa"
$ws.Range("D207").Value = "ASM"

# Bermuda / BMU  (name typed before code for this row)
$ws.Range("B208").Value = "Bermuda"
$ws.Range("C208").Value = "Bermuda"
$ws.Range("A208").Value = "BMU"
$ws.Range("D208").Value = "BMU"

# GRL / Greenland
$ws.Range("A209").Value = "GRL"
$ws.Range("B209").Value = "Greenland"
$ws.Range("C209").Value = "Greenland"
$ws.Range("D209").Value = "GRL"

# PRI / Puerto Rico
$ws.Range("A210").Value = "PRI"
$ws.Range("B210").Value = "Puerto Rico"
$ws.Range("C210").Value = "Puerto Rico"
$ws.Range("D210").Value = "PRI"

# PYF / French Polynesia
$ws.Range("A211").Value = "PYF"
$ws.Range("B211").Value = "French Polynesia"
$ws.Range("C211").Value = "French Polynesia"
$ws.Range("D211").Value = "PYF"

$newRange = $ws.Range("A207:D211")
$newRange.Interior.Color = 65535

$ws.Range("A194").Select()
$ws.Application.ActiveWindow.ScrollRow = 194
$ws.Range("E202").Select()
